# "save data done + era data updated"
# Adds a new "Save" column (H) to Sheet1: H1 gets the header label "Save"
# (styled like the other header cells), and H2:H51 get a 0/1 flag derived
# from the existing "sum" column (G) — 1 when the sum is large (>= 9),
# else 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 -------------------------------------------------
# Clone the formatting of the neighbouring header cell (G1: bold text,
# thin border, centered) onto H1, then set its own label.
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(1, 8))
$ws.Cells.Item(1, 8).Value = "Save"

# --- Data rows 2-51 ---------------------------------------------------
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $sumVal = $ws.Cells.Item($r, 7).Value2
    if ($sumVal -ge 9) {
        $flag = 1
    } else {
        $flag = 0
    }
    $ws.Cells.Item($r, 8).Value = $flag
}

Write-Output "Save column (H) written for rows 1-$lastRow"
